$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.60954197817042
$ws.Range("C2").Value = 6.393399622049882
$ws.Range("D2").Value = 6.484671421832469
$ws.Range("E2").Value = 16.30822471466244
$ws.Range("F2").Value = 35.06749170459285
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("K2").Value = 12.91713939145584
$ws.Range("N2").Value = 20.89185839143319
$ws.Range("B3").Value = 13.29928857578269
$ws.Range("C3").Value = 6.111244584545926
$ws.Range("D3").Value = 6.498699826115866
$ws.Range("E3").Value = 15.39321256051893
$ws.Range("F3").Value = 34.78736585981625
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("K3").Value = 12.69698921044773
$ws.Range("N3").Value = 20.91959389422701
$ws.Range("B4").Value = 13.10965990068858
$ws.Range("C4").Value = 5.93345404340446
$ws.Range("D4").Value = 6.507556599172942
$ws.Range("E4").Value = 14.80825640629554
$ws.Range("F4").Value = 34.62315099477182
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("K4").Value = 12.56406682283925
$ws.Range("N4").Value = 20.93845682681559
$ws.Range("B5").Value = 13.03272826527347
$ws.Range("C5").Value = 5.859983051374587
$ws.Range("D5").Value = 6.511227304034957
$ws.Range("E5").Value = 14.56434032838408
$ws.Range("F5").Value = 34.55824208636061
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("K5").Value = 12.51054216047668
$ws.Range("N5").Value = 20.94660367206159
$ws.Range("B6").Value = 13.0199782856309
$ws.Range("C6").Value = 5.847725551560162
$ws.Range("D6").Value = 6.511840545584141
$ws.Range("E6").Value = 14.52351264597437
$ws.Range("F6").Value = 34.54758687092081
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("K6").Value = 12.50169543769643
$ws.Range("N6").Value = 20.94798420787851
$ws.Range("B7").Value = 13.1086208153055
$ws.Range("C7").Value = 5.932467143345951
$ws.Range("D7").Value = 6.507605854197755
$ws.Range("E7").Value = 14.80498890344661
$ws.Range("F7").Value = 34.6222674057789
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("K7").Value = 12.56334226974253
$ws.Range("N7").Value = 20.93856483664765
$ws.Range("B8").Value = 13.50246073389677
$ws.Range("C8").Value = 6.29712544017027
$ws.Range("D8").Value = 6.489458120677437
$ws.Range("E8").Value = 15.99768001179853
$ws.Range("F8").Value = 34.96932169457136
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("K8").Value = 12.84081046580285
$ws.Range("N8").Value = 20.90104054740191
$ws.Range("B9").Value = 14.27632484620022
$ws.Range("C9").Value = 6.971478470939204
$ws.Range("D9").Value = 6.455785322062637
$ws.Range("E9").Value = 18.16694151412842
$ws.Range("F9").Value = 35.70915304109386
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("K9").Value = 13.39941147107741
$ws.Range("N9").Value = 20.84204368229384
$ws.Range("B10").Value = 14.83898125040868
$ws.Range("C10").Value = 7.436798923952876
$ws.Range("D10").Value = 6.4321907441379
$ws.Range("E10").Value = 19.79230145112534
$ws.Range("F10").Value = 36.28525387245356
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("K10").Value = 13.81416497874525
$ws.Range("N10").Value = 20.80765111702445
$ws.Range("B11").Value = 15.09232927512537
$ws.Range("C11").Value = 7.641072030997474
$ws.Range("D11").Value = 6.421700686622354
$ws.Range("E11").Value = 20.49057219916619
$ws.Range("F11").Value = 36.55358423532585
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("K11").Value = 14.00286883827805
$ws.Range("N11").Value = 20.79396176333289
$ws.Range("B12").Value = 15.18778806604171
$ws.Range("C12").Value = 7.717301828672779
$ws.Range("D12").Value = 6.417762997475651
$ws.Range("E12").Value = 20.74911539155547
$ws.Range("F12").Value = 36.6560218819017
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("K12").Value = 14.07425658488964
$ws.Range("N12").Value = 20.78906027981356
$ws.Range("B13").Value = 15.16725216986968
$ws.Range("C13").Value = 7.700935216203633
$ws.Range("D13").Value = 6.418609512301609
$ws.Range("E13").Value = 20.69369430052288
$ws.Range("F13").Value = 36.63392451363512
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("K13").Value = 14.05888624086158
$ws.Range("N13").Value = 20.7901033255147
$ws.Range("B14").Value = 15.10019288827483
$ws.Range("C14").Value = 7.647366336390227
$ws.Range("D14").Value = 6.421376037533133
$ws.Range("E14").Value = 20.51196028111629
$ws.Range("F14").Value = 36.56199569761709
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("K14").Value = 14.00874377750355
$ws.Range("N14").Value = 20.79355284973218
$ws.Range("B15").Value = 15.05905185548993
$ws.Range("C15").Value = 7.614405868361492
$ws.Range("D15").Value = 6.423075119774717
$ws.Range("E15").Value = 20.39987866069176
$ws.Range("F15").Value = 36.51804268920964
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("K15").Value = 13.97801874830455
$ws.Range("N15").Value = 20.7957025896553
$ws.Range("B16").Value = 14.82236287387841
$ws.Range("C16").Value = 7.423294929176697
$ws.Range("D16").Value = 6.432881162112835
$ws.Range("E16").Value = 19.74584352617361
$ws.Range("F16").Value = 36.26783757199052
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("K16").Value = 13.8018268995754
$ws.Range("N16").Value = 20.80858520284805
$ws.Range("B17").Value = 14.67641809502877
$ws.Range("C17").Value = 7.304113088791842
$ws.Range("D17").Value = 6.438958929172722
$ws.Range("E17").Value = 19.33410914472458
$ws.Range("F17").Value = 36.11589631629958
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("K17").Value = 13.69369265663082
$ws.Range("N17").Value = 20.81699003489425
$ws.Range("B18").Value = 14.592234884122
$ws.Range("C18").Value = 7.234868908005339
$ws.Range("D18").Value = 6.442477612752498
$ws.Range("E18").Value = 19.09341850779617
$ws.Range("F18").Value = 36.02909737193071
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("K18").Value = 13.63150386155309
$ws.Range("N18").Value = 20.82200831358856
$ws.Range("B19").Value = 14.56369401187707
$ws.Range("C19").Value = 7.21130679703285
$ws.Range("D19").Value = 6.443672923258386
$ws.Range("E19").Value = 19.01125871118485
$ws.Range("F19").Value = 35.99981290815454
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("K19").Value = 13.61045135656502
$ws.Range("N19").Value = 20.82373899202278
$ws.Range("B20").Value = 14.69197974060214
$ws.Range("C20").Value = 7.316872478736756
$ws.Range("D20").Value = 6.438309571436641
$ws.Range("E20").Value = 19.37833944310921
$ws.Range("F20").Value = 36.13200980873501
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("K20").Value = 13.70520354053465
$ws.Range("N20").Value = 20.81607627072473
$ws.Range("B21").Value = 15.11990360562714
$ws.Range("C21").Value = 7.663131747590468
$ws.Range("D21").Value = 6.420562503580269
$ws.Range("E21").Value = 20.565499181705
$ws.Range("F21").Value = 36.58310108768828
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("K21").Value = 14.0234743231299
$ws.Range("N21").Value = 20.7925319683769
$ws.Range("B22").Value = 15.39674379772434
$ws.Range("C22").Value = 7.882859397866519
$ws.Range("D22").Value = 6.409165698469921
$ws.Range("E22").Value = 21.30713855461097
$ws.Range("F22").Value = 36.88269940358288
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("K22").Value = 14.2310403880349
$ws.Range("N22").Value = 20.77879083166575
$ws.Range("B23").Value = 15.2492806101463
$ws.Range("C23").Value = 7.766205323707675
$ws.Range("D23").Value = 6.415230012739907
$ws.Range("E23").Value = 20.91443283722312
$ws.Range("F23").Value = 36.72238496427234
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("K23").Value = 14.12032283203509
$ws.Range("N23").Value = 20.78597372117633
$ws.Range("B24").Value = 14.68494517735681
$ws.Range("C24").Value = 7.311106212039063
$ws.Range("D24").Value = 6.438603069710254
$ws.Range("E24").Value = 19.35835532231437
$ws.Range("F24").Value = 36.12472316713596
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("K24").Value = 13.69999952952019
$ws.Range("N24").Value = 20.81648880343234
$ws.Range("B25").Value = 14.06753569357277
$ws.Range("C25").Value = 6.794005497978477
$ws.Range("D25").Value = 6.464692042447573
$ws.Range("E25").Value = 17.58471719701074
$ws.Range("F25").Value = 35.50302250119643
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("K25").Value = 13.24721134675389
$ws.Range("N25").Value = 20.85643602709043
